# Optuna Attempt (go back with original)
# Updates forecast "Inventory Coverage" (H) and "Seasonality Index" (L)
# values on the "Forecast Comparison" sheet, and the 16/8/4-week total
# forecast figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Forecast Comparison ------------------------------------------
$fc = $wb.Worksheets.Item("Forecast Comparison")

# row : H (Inventory Coverage) , L (Seasonality Index)
$fc.Cells.Item(2, 8).Value  = 20.17
$fc.Cells.Item(2, 12).Value = 1.08

$fc.Cells.Item(3, 8).Value  = 19.17
$fc.Cells.Item(3, 12).Value = 0.99

$fc.Cells.Item(4, 8).Value  = 16.38
$fc.Cells.Item(4, 12).Value = 1.09

$fc.Cells.Item(5, 8).Value  = 15.38
$fc.Cells.Item(5, 12).Value = 0.84

$fc.Cells.Item(6, 8).Value  = 14.38
$fc.Cells.Item(6, 12).Value = 0.82

$fc.Cells.Item(7, 8).Value  = 14.84
$fc.Cells.Item(7, 12).Value = 1.09

$fc.Cells.Item(8, 8).Value  = 12.48
$fc.Cells.Item(8, 12).Value = 0.91

$fc.Cells.Item(9, 8).Value  = 11.48
$fc.Cells.Item(9, 12).Value = 1.17

$fc.Cells.Item(10, 8).Value  = 10.48
$fc.Cells.Item(10, 12).Value = 1.2

$fc.Cells.Item(11, 8).Value  = 9.48
$fc.Cells.Item(11, 12).Value = 1.03

$fc.Cells.Item(12, 8).Value  = 8.48
$fc.Cells.Item(12, 12).Value = 0.92

$fc.Cells.Item(13, 8).Value  = 7.48
$fc.Cells.Item(13, 12).Value = 1.16

$fc.Cells.Item(14, 8).Value  = 6.48
$fc.Cells.Item(14, 12).Value = 0.94

$fc.Cells.Item(15, 8).Value  = 6.08
$fc.Cells.Item(15, 12).Value = 1.1

$fc.Cells.Item(16, 8).Value  = 5.49
$fc.Cells.Item(16, 12).Value = 1.04

$fc.Cells.Item(17, 8).Value  = 4.15
$fc.Cells.Item(17, 12).Value = 1.14

# ---- Sheet: Summary --------------------------------------------------------
# B9/B10/B11 hold the forecast totals as text labels ("40", "20", "10") --
# format the cells as Text first so the numeric-looking strings are not
# auto-converted to numbers, then restore the default "Normal" style so no
# extra formatting is left behind on the cells.
$sm = $wb.Worksheets.Item("Summary")

$sm.Range("B9").NumberFormat = "@"
$sm.Range("B9").Value = "40"
$sm.Range("B9").Style = "Normal"

$sm.Range("B10").NumberFormat = "@"
$sm.Range("B10").Value = "20"
$sm.Range("B10").Style = "Normal"

$sm.Range("B11").NumberFormat = "@"
$sm.Range("B11").Value = "10"
$sm.Range("B11").Style = "Normal"

Write-Output "Applied forecast summary updates (Optuna attempt revert)."
